$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that gets bumped by one day
# for every data row (rows 2-29) as part of the automatic update.
for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $cell.Value2 + 1
}
